$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel's automatic
# "looks like a number" type coercion (which would otherwise turn strings
# like "591.47" into real numbers, or force a quote-prefixed Text style).
# We do this by entering a formula that evaluates to the literal string,
# then collapsing it to a plain stored value via Copy + PasteSpecial
# (xlPasteValues = -4163). This keeps the cell's style/number-format
# untouched, matching the original inline-string cells.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Formula = '="' + $val + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue 'D2' '68.667.97'
Set-TextValue 'E2' '  +1.33%  '
Set-TextValue 'D3' '2.508.03'
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '591.47'
Set-TextValue 'E5' '  +0.72%  '
Set-TextValue 'D6' '174.04'
Set-TextValue 'E6' '  -1.10%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'D8' '0.515'
Set-TextValue 'E8' '  -0.29%  '
Set-TextValue 'D9' '2.507.23'
Set-TextValue 'E9' '  +0.35%  '
Set-TextValue 'E10' '  +5.20%  '
Set-TextValue 'E11' '  -1.37%  '
Set-TextValue 'D12' '5.00'
Set-TextValue 'E12' '  +1.06%  '
Set-TextValue 'D13' '0.334'
Set-TextValue 'E13' '  -1.58%  '
Set-TextValue 'D14' '2.974.67'
Set-TextValue 'E14' '  +1.05%  '
Set-TextValue 'D15' '25.68'
Set-TextValue 'E15' '  -0.39%  '
Set-TextValue 'D16' '68.535.98'
Set-TextValue 'E16' '  +1.38%  '
Set-TextValue 'E17' '  -0.82%  '
Set-TextValue 'D18' '2.513.39'
Set-TextValue 'E18' '  +5.21%  '
Set-TextValue 'D19' '360.91'
Set-TextValue 'E19' '  +2.40%  '
Set-TextValue 'E20' '  +0.40%  '
Set-TextValue 'D21' '10.84'
Set-TextValue 'E21' '  -2.37%  '
Set-TextValue 'D22' '4.00'
Set-TextValue 'E22' '  -2.08%  '
Set-TextValue 'E23' '  -0.07%  '
Set-TextValue 'D24' '70.04'
Set-TextValue 'E24' '  -1.01%  '
Set-TextValue 'E25' '  -2.56%  '
Set-TextValue 'E26' '  -3.81%  '
Set-TextValue 'D27' '1.65'
Set-TextValue 'E27' '  -7.32%  '
Set-TextValue 'D28' '2.638.79'
Set-TextValue 'E28' '  +0.57%  '
Set-TextValue 'E29' '  +0.37%  '
Set-TextValue 'D30' '507.21'
Set-TextValue 'E30' '  -1.40%  '
Set-TextValue 'D31' '0.0₃0872'
Set-TextValue 'E31' '  -4.58%  '
Set-TextValue 'D32' '7.70'
Set-TextValue 'E32' '  -1.98%  '
Set-TextValue 'E33' '  -3.06%  '
Set-TextValue 'E34' '  -0.83%  '
Set-TextValue 'E35' '  +0.07%  '
Set-TextValue 'D36' '161.38'
Set-TextValue 'E36' '  -0.29%  '
Set-TextValue 'E37' '  -4.23%  '
Set-TextValue 'D38' '18.52'
Set-TextValue 'E38' '  +0.19%  '
Set-TextValue 'E39' '  -0.28%  '
Set-TextValue 'E41' '  -2.54%  '
Set-TextValue 'D42' '1.69'
Set-TextValue 'E42' '  -4.29%  '
Set-TextValue 'D43' '4.72'
Set-TextValue 'E43' '  -3.15%  '
Set-TextValue 'D44' '0.316'
Set-TextValue 'E44' '  -4.49%  '
Set-TextValue 'D45' '2.31'
Set-TextValue 'E45' '  -5.38%  '
Set-TextValue 'D46' '149.75'
Set-TextValue 'E46' '  +3.13%  '
Set-TextValue 'E47' '  +0.34%  '
Set-TextValue 'D48' '0.511'
Set-TextValue 'E48' '  -0.93%  '
Set-TextValue 'D49' '0.0736'
Set-TextValue 'E49' '  -1.42%  '
Set-TextValue 'D50' '0.0₆0247'
Set-TextValue 'E50' '  -4.46%  '
Set-TextValue 'E51' '  -2.49%  '
